$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 338, shifting existing rows 338:411 down to 339:412
$ws.Rows("338:338").Insert()

# Populate the new record in row 338
$ws.Range("A338").Value = 4
$ws.Range("B338").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C338").Value = "Los Lagos"
$ws.Range("D338").Value = 45015
$ws.Range("E338").Value = 10
$ws.Range("F338").Value = 100112043
$ws.Range("G338").Value = "Pepino ensalada"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 200
$ws.Range("K338").Value = 14000
$ws.Range("L338").Value = 14000
$ws.Range("M338").Value = 14000
$ws.Range("N338").Value = "$/caja 60 unidades"
$ws.Range("O338").Value = "Región de Arica y Parinacota"
$ws.Range("P338").Value = 233
$ws.Range("Q338").Value = 60
$ws.Range("R338").Value = "Hortaliza"
